$d = $word.ActiveDocument

# Remove the (hidden) "_GoBack" bookmark that wraps the end of the first
# paragraph - it carries no visible text, so deleting it just drops the
# <w:bookmarkStart/>/<w:bookmarkEnd/> markers.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Add a blank paragraph, then a new paragraph containing the new line of
# text, both appended after the existing "Que coisa feia!" paragraph.
$tail = $d.Range($d.Content.End, $d.Content.End)
$tail.InsertParagraphAfter()

$tail2 = $d.Range($d.Content.End, $d.Content.End)
$tail2.InsertParagraphAfter()

$tail3 = $d.Range($d.Content.End, $d.Content.End)
$tail3.Text = "Feia é meus bago no frio"
